$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("run_1")
$ws.Cells.Item(2, 6).Value = 27.11155676841736
$ws.Cells.Item(3, 6).Value = 26.71789813041687
$ws.Cells.Item(4, 6).Value = 26.61899018287659
$ws.Cells.Item(5, 6).Value = 26.70770621299744
$ws.Cells.Item(6, 6).Value = 26.71150660514832
$ws.Cells.Item(7, 6).Value = 26.75455546379089
$ws.Cells.Item(8, 6).Value = 26.7333972454071
$ws.Cells.Item(9, 6).Value = 26.74641966819763
$ws.Cells.Item(10, 6).Value = 26.83454918861389
$ws.Cells.Item(11, 6).Value = 26.91780710220337
$ws.Cells.Item(12, 6).Value = 26.70109581947327
$ws.Cells.Item(13, 6).Value = 26.65547919273376
$ws.Cells.Item(14, 6).Value = 26.9501051902771
$ws.Cells.Item(15, 6).Value = 26.64055585861206
$ws.Cells.Item(16, 6).Value = 26.66196084022522
$ws.Cells.Item(17, 6).Value = 26.78947281837464
$ws.Cells.Item(18, 6).Value = 26.68064761161804
$ws.Cells.Item(19, 6).Value = 26.76467680931092
$ws.Cells.Item(20, 6).Value = 26.69610357284546
$ws.Cells.Item(21, 6).Value = 27.07684278488159

$ws = $wb.Worksheets.Item("run_2")
$ws.Cells.Item(2, 6).Value = 26.92121696472168
$ws.Cells.Item(3, 6).Value = 26.79382705688477
$ws.Cells.Item(4, 6).Value = 26.71684694290161
$ws.Cells.Item(5, 6).Value = 26.77211308479309
$ws.Cells.Item(6, 6).Value = 26.79521560668945
$ws.Cells.Item(7, 6).Value = 26.69760656356812
$ws.Cells.Item(8, 6).Value = 26.84830617904663
$ws.Cells.Item(9, 6).Value = 26.69267249107361
$ws.Cells.Item(10, 6).Value = 26.99621367454529
$ws.Cells.Item(11, 6).Value = 27.0455482006073
$ws.Cells.Item(12, 6).Value = 26.8530387878418
$ws.Cells.Item(13, 6).Value = 26.89045310020447
$ws.Cells.Item(14, 6).Value = 26.65666818618774
$ws.Cells.Item(15, 6).Value = 26.77835559844971
$ws.Cells.Item(16, 6).Value = 26.68855929374695
$ws.Cells.Item(17, 6).Value = 26.74869871139526
$ws.Cells.Item(18, 6).Value = 26.73119115829468
$ws.Cells.Item(19, 6).Value = 26.73916411399841
$ws.Cells.Item(20, 6).Value = 26.81686234474182
$ws.Cells.Item(21, 6).Value = 26.92924523353577

$ws = $wb.Worksheets.Item("run_3")
$ws.Cells.Item(2, 6).Value = 27.01069140434265
$ws.Cells.Item(3, 6).Value = 26.69698977470398
$ws.Cells.Item(4, 6).Value = 26.71156001091004
$ws.Cells.Item(5, 6).Value = 26.71007394790649
$ws.Cells.Item(6, 6).Value = 26.743967294693
$ws.Cells.Item(7, 6).Value = 26.74365186691284
$ws.Cells.Item(8, 6).Value = 26.81051468849182
$ws.Cells.Item(9, 6).Value = 26.64302849769592
$ws.Cells.Item(10, 6).Value = 26.70071911811829
$ws.Cells.Item(11, 6).Value = 26.94250249862671
$ws.Cells.Item(12, 6).Value = 26.715411901474
$ws.Cells.Item(13, 6).Value = 26.73483943939209
$ws.Cells.Item(14, 6).Value = 26.67713260650635
$ws.Cells.Item(15, 6).Value = 26.7915632724762
$ws.Cells.Item(16, 6).Value = 26.80637645721436
$ws.Cells.Item(17, 6).Value = 26.75828838348389
$ws.Cells.Item(18, 6).Value = 26.77852940559387
$ws.Cells.Item(19, 6).Value = 26.70512557029724
$ws.Cells.Item(20, 6).Value = 26.64127516746521
$ws.Cells.Item(21, 6).Value = 26.97547578811645

$ws = $wb.Worksheets.Item("run_4")
$ws.Cells.Item(2, 6).Value = 27.10856127738953
$ws.Cells.Item(3, 6).Value = 26.76086354255676
$ws.Cells.Item(4, 6).Value = 26.73359155654907
$ws.Cells.Item(5, 6).Value = 26.70680046081543
$ws.Cells.Item(6, 6).Value = 26.75707626342773
$ws.Cells.Item(7, 6).Value = 26.68600130081177
$ws.Cells.Item(8, 6).Value = 26.75049138069153
$ws.Cells.Item(9, 6).Value = 26.75858068466187
$ws.Cells.Item(10, 6).Value = 26.85736227035522
$ws.Cells.Item(11, 6).Value = 26.95536780357361
$ws.Cells.Item(12, 6).Value = 26.710857629776
$ws.Cells.Item(13, 6).Value = 26.88304400444031
$ws.Cells.Item(14, 6).Value = 26.89371085166931
$ws.Cells.Item(15, 6).Value = 26.86560773849488
$ws.Cells.Item(16, 6).Value = 26.74056816101075
$ws.Cells.Item(17, 6).Value = 26.9486255645752
$ws.Cells.Item(18, 6).Value = 26.81576442718506
$ws.Cells.Item(19, 6).Value = 26.89855408668518
$ws.Cells.Item(20, 6).Value = 26.81344199180603
$ws.Cells.Item(21, 6).Value = 26.89301753044128

$ws = $wb.Worksheets.Item("run_5")
$ws.Cells.Item(2, 6).Value = 27.05046677589417
$ws.Cells.Item(3, 6).Value = 26.85525345802307
$ws.Cells.Item(4, 6).Value = 26.87679696083069
$ws.Cells.Item(5, 6).Value = 26.96823835372925
$ws.Cells.Item(6, 6).Value = 26.76120114326477
$ws.Cells.Item(7, 6).Value = 26.69966125488281
$ws.Cells.Item(8, 6).Value = 26.78378415107727
$ws.Cells.Item(9, 6).Value = 26.8333580493927
$ws.Cells.Item(10, 6).Value = 26.67384791374207
$ws.Cells.Item(11, 6).Value = 26.96737027168274
$ws.Cells.Item(12, 6).Value = 26.90947723388672
$ws.Cells.Item(13, 6).Value = 26.84405064582825
$ws.Cells.Item(14, 6).Value = 26.85953330993652
$ws.Cells.Item(15, 6).Value = 26.91600227355957
$ws.Cells.Item(16, 6).Value = 26.75749969482422
$ws.Cells.Item(17, 6).Value = 26.67129802703857
$ws.Cells.Item(18, 6).Value = 26.67506551742554
$ws.Cells.Item(19, 6).Value = 26.72724390029907
$ws.Cells.Item(20, 6).Value = 26.8412299156189
$ws.Cells.Item(21, 6).Value = 26.88141226768494
